# Updates from Catalina - Added some ship and vessel info
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 60-63: Madison Swanson / Twin Ridges entries
$ws.Range("O60").Value = "Kongsberg Simrad EM3000"
$ws.Range("P60").Value = 300

$ws.Range("O61").Value = "Kongsberg Simrad EM3000"
$ws.Range("P61").Value = 300

$ws.Range("N62").Value = "R/V Suncoaster"
$ws.Range("O62").Value = "Kongsberg Simrad EM3000"
$ws.Range("P62").Value = 300

$ws.Range("N63").Value = "R/V Suncoaster"
$ws.Range("O63").Value = "Kongsberg Simrad EM3000"
$ws.Range("P63").Value = 300

# Rows 104-109: DeSoto Canyon entries
$ws.Range("N104").Value = "R/V Moana Wave"
$ws.Range("O104").Value = "Kongsberg Dual-head EM3000"
$ws.Range("P104").Value = 300

$ws.Range("N105").Value = "R/V Moana Wave"
$ws.Range("O105").Value = "Kongsberg Dual-head EM3000"
$ws.Range("P105").Value = 300

$ws.Range("N106").Value = "R/V Moana Wave"
$ws.Range("O106").Value = "Kongsberg Dual-head EM3000"
$ws.Range("P106").Value = 300

$ws.Range("N107").Value = "R/V Moana Wave"
$ws.Range("O107").Value = "Kongsberg Dual-head EM3000"
$ws.Range("P107").Value = 300

$ws.Range("N108").Value = "R/V Moana Wave"
$ws.Range("O108").Value = "Kongsberg Dual-head EM3000"
$ws.Range("P108").Value = 300

$ws.Range("N109").Value = "R/V Moana Wave"
$ws.Range("O109").Value = "Kongsberg Dual-head EM3000"
$ws.Range("P109").Value = 300

# Rows 112-113: Pulley Ridge entries
$ws.Range("N112").Value = "R/V's Bellows & Suncoaster"
$ws.Range("O112").Value = "Kongsberg Simrad EM3000"
$ws.Range("P112").Value = 300

$ws.Range("N113").Value = "R/V's Bellows & Suncoaster"
$ws.Range("O113").Value = "Kongsberg Simrad EM3000"
$ws.Range("P113").Value = 300

# Rows 114-121: Steamboat Lumps / West Florida Shelf entries
$ws.Range("N114").Value = "R/V Moana Wave"
$ws.Range("O114").Value = "Kongsberg EM1002"
$ws.Range("P114").Value = 95

$ws.Range("N115").Value = "R/V Moana Wave"
$ws.Range("O115").Value = "Kongsberg EM1002"
$ws.Range("P115").Value = 95

$ws.Range("N116").Value = "R/V Moana Wave"
$ws.Range("O116").Value = "Kongsberg EM1002"
$ws.Range("P116").Value = 95

$ws.Range("N117").Value = "R/V Moana Wave"
$ws.Range("O117").Value = "Kongsberg EM1002"
$ws.Range("P117").Value = 95

$ws.Range("N118").Value = "R/V Moana Wave"
$ws.Range("O118").Value = "Kongsberg EM1002"
$ws.Range("P118").Value = 95

$ws.Range("N119").Value = "R/V Moana Wave"
$ws.Range("O119").Value = "Kongsberg EM1002"
$ws.Range("P119").Value = 95

$ws.Range("N120").Value = "R/V Moana Wave"
$ws.Range("O120").Value = "Kongsberg EM1002"
$ws.Range("P120").Value = 95

$ws.Range("N121").Value = "R/V Moana Wave"
$ws.Range("O121").Value = "Kongsberg EM1002"
$ws.Range("P121").Value = 95
